$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = [double]"0.922019198695378"
$ws.Range("C2").Value = [double]"0.00287239168303621"
$ws.Range("D2").Value = [double]"0.00196434528001186"
$ws.Range("E2").Value = [double]"0.990011489566732"
$ws.Range("F2").Value = [double]"0.0227196916348542"
$ws.Range("G2").Value = [double]"0.00389162744153293"
$ws.Range("H2").Value = [double]"0.948760238686483"
$ws.Range("I2").Value = [double]"0.994866758089026"
$ws.Range("J2").Value = [double]"0.000426225862644083"
$ws.Range("K2").Value = [double]"0.981524035432341"
$ws.Range("L2").Value = [double]"0.0244801897631667"
$ws.Range("M2").Value = [double]"0.000333568066417108"
$ws.Range("N2").Value = [double]"0.948778770245728"
$ws.Range("O2").Value = [double]"0.00168637189133094"
$ws.Range("P2").Value = [double]"0.000148252473963159"
$ws.Range("Q2").Value = [double]"0.132760090434009"
$ws.Range("R2").Value = [double]"0.00164930877284015"
$ws.Range("S2").Value = [double]"0.999332863867166"
$ws.Range("T2").Value = [double]"0.000148252473963159"
$ws.Range("U2").Value = [double]"0.00402134835625069"
$ws.Range("V2").Value = [double]"0.000259441829435529"
$ws.Range("W2").Value = [double]"0.0532782328305104"
$ws.Range("X2").Value = [double]"0.000481820540380268"
$ws.Range("B3").Value = [double]"0.000537415218116452"
$ws.Range("C3").Value = [double]"0.000352099625662503"
$ws.Range("D3").Value = [double]"0.994885289648271"
$ws.Range("E3").Value = [double]"0.000796857047551981"
$ws.Range("F3").Value = [double]"0.000129720914717764"
$ws.Range("G3").Value = [double]"0.00374337496756977"
$ws.Range("H3").Value = [double]"0.00211259775397502"
$ws.Range("I3").Value = [double]"0.00207553463548423"
$ws.Range("J3").Value = [double]"0.994014306363737"
$ws.Range("K3").Value = [double]"0.00170490345057633"
$ws.Range("L3").Value = [double]"1.85315592453949e-05"
$ws.Range("M3").Value = [double]"0.000500352099625662"
$ws.Range("N3").Value = [double]"0.0047255476075757"
$ws.Range("O3").Value = [double]"0.913124050257589"
$ws.Range("P3").Value = [double]"0.000593009895852637"
$ws.Range("Q3").Value = [double]"0.000389162744153293"
$ws.Range("R3").Value = [double]"0.00552240465512768"
$ws.Range("S3").Value = [double]"3.70631184907898e-05"
$ws.Range("T3").Value = [double]"0.00748674993513954"
$ws.Range("U3").Value = [double]"1.85315592453949e-05"
$ws.Range("V3").Value = [double]"0.00151958785812238"
$ws.Range("W3").Value = [double]"0.00187168748378489"
$ws.Range("X3").Value = [double]"0.00709758719098625"
$ws.Range("B4").Value = [double]"0.055186983432786"
$ws.Range("C4").Value = [double]"0.00114895667321448"
$ws.Range("D4").Value = [double]"0.000333568066417108"
$ws.Range("E4").Value = [double]"0.00478114228531189"
$ws.Range("F4").Value = [double]"0.973333086245877"
$ws.Range("G4").Value = [double]"0.00559653089210926"
$ws.Range("H4").Value = [double]"0.0408620881360958"
$ws.Range("I4").Value = [double]"0.000870983284533561"
$ws.Range("J4").Value = [double]"0.00498498943701123"
$ws.Range("K4").Value = [double]"0.00674548756532375"
$ws.Range("L4").Value = [double]"0.970164189614914"
$ws.Range("M4").Value = [double]"0.00157518253585857"
$ws.Range("N4").Value = [double]"0.0317445609873615"
$ws.Range("O4").Value = [double]"0.000277973388680924"
$ws.Range("P4").Value = [double]"0"
$ws.Range("Q4").Value = [double]"0.846966383751529"
$ws.Range("R4").Value = [double]"0.992680034098069"
$ws.Range("S4").Value = [double]"0.000481820540380268"
$ws.Range("T4").Value = [double]"5.55946777361847e-05"
$ws.Range("U4").Value = [double]"0.989529669026352"
$ws.Range("V4").Value = [double]"0.00413253771172306"
$ws.Range("W4").Value = [double]"0.938771728253215"
$ws.Range("X4").Value = [double]"0.99227233979467"
$ws.Range("B5").Value = [double]"0.0147511211593343"
$ws.Range("C5").Value = [double]"0.995181794596197"
$ws.Range("D5").Value = [double]"0.00248322893888292"
$ws.Range("E5").Value = [double]"0.00283532856454542"
$ws.Range("F5").Value = [double]"0.000315036507171713"
$ws.Range("G5").Value = [double]"0.985211815722175"
$ws.Range("H5").Value = [double]"0.00153811941736778"
$ws.Range("I5").Value = [double]"0.00159371409510396"
$ws.Range("J5").Value = [double]"0.000333568066417108"
$ws.Range("K5").Value = [double]"0.00743115525740336"
$ws.Range("L5").Value = [double]"5.55946777361847e-05"
$ws.Range("M5").Value = [double]"0.997498239501872"
$ws.Range("N5").Value = [double]"0.00776472332382047"
$ws.Range("O5").Value = [double]"0.0767391868351803"
$ws.Range("P5").Value = [double]"0.999166079833957"
$ws.Range("Q5").Value = [double]"0.00322449130869871"
$ws.Range("R5").Value = [double]"1.85315592453949e-05"
$ws.Range("S5").Value = [double]"0"
$ws.Range("T5").Value = [double]"0.990104147362959"
$ws.Range("U5").Value = [double]"0.00557799933286387"
$ws.Range("V5").Value = [double]"0.993773396093547"
$ws.Range("W5").Value = [double]"0.000685667692079612"
$ws.Range("X5").Value = [double]"5.55946777361847e-05"
